# Update cryptos list (Price + Volume(1h) columns) per upstream scrape refresh.
# Column D ("Price") values are written as literal text via an explicit Text
# number format (then ClearFormats restores the default/unstyled cell) so that
# numeric-looking strings like "289.80" are not auto-converted to the number 289.8
# by Excel's normal cell-entry parsing -- matching the source data, which uses
# dotted thousand-separator price strings (e.g. "22.420.89") throughout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = '@'
$dCell.Value = '22.420.89'
$dCell.ClearFormats()
$ws.Range("E2").Value = '  -3.77%  '
$dCell = $ws.Range("D3")
$dCell.NumberFormat = '@'
$dCell.Value = '1.572.24'
$dCell.ClearFormats()
$ws.Range("E3").Value = '  -3.46%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("E5").Value = '  -0.10%  '
$dCell = $ws.Range("D6")
$dCell.NumberFormat = '@'
$dCell.Value = '289.80'
$dCell.ClearFormats()
$ws.Range("E6").Value = '  -2.74%  '
$dCell = $ws.Range("D7")
$dCell.NumberFormat = '@'
$dCell.Value = '0.3670'
$dCell.ClearFormats()
$ws.Range("E7").Value = '  -2.52%  '
$dCell = $ws.Range("D8")
$dCell.NumberFormat = '@'
$dCell.Value = '49.28'
$dCell.ClearFormats()
$ws.Range("E8").Value = '  -1.56%  '
$dCell = $ws.Range("D9")
$dCell.NumberFormat = '@'
$dCell.Value = '0.3403'
$dCell.ClearFormats()
$ws.Range("E9").Value = '  -2.00%  '
$dCell = $ws.Range("D10")
$dCell.NumberFormat = '@'
$dCell.Value = '1.172'
$dCell.ClearFormats()
$ws.Range("E10").Value = '  -2.16%  '
$dCell = $ws.Range("D11")
$dCell.NumberFormat = '@'
$dCell.Value = '0.07651'
$dCell.ClearFormats()
$ws.Range("E11").Value = '  -4.66%  '
$dCell = $ws.Range("D12")
$dCell.NumberFormat = '@'
$dCell.Value = '1.001'
$dCell.ClearFormats()
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("E13").Value = '  -2.70%  '
$dCell = $ws.Range("D14")
$dCell.NumberFormat = '@'
$dCell.Value = '6.063'
$dCell.ClearFormats()
$ws.Range("E14").Value = '  -3.52%  '
$dCell = $ws.Range("D15")
$dCell.NumberFormat = '@'
$dCell.Value = '6.924'
$dCell.ClearFormats()
$ws.Range("E15").Value = '  -3.97%  '
$dCell = $ws.Range("D16")
$dCell.NumberFormat = '@'
$dCell.Value = '1.569.90'
$dCell.ClearFormats()
$ws.Range("E16").Value = '  -3.37%  '
$dCell = $ws.Range("D17")
$dCell.NumberFormat = '@'
$dCell.Value = '0.00001133'
$dCell.ClearFormats()
$ws.Range("E17").Value = '  -4.60%  '
$dCell = $ws.Range("D18")
$dCell.NumberFormat = '@'
$dCell.Value = '90.02'
$dCell.ClearFormats()
$ws.Range("E18").Value = '  -4.79%  '
$dCell = $ws.Range("D19")
$dCell.NumberFormat = '@'
$dCell.Value = '0.06739'
$dCell.ClearFormats()
$ws.Range("E19").Value = '  -2.92%  '
$dCell = $ws.Range("D21")
$dCell.NumberFormat = '@'
$dCell.Value = '6.265'
$dCell.ClearFormats()
$ws.Range("E21").Value = '  -5.07%  '
$ws.Range("E22").Value = '  -3.72%  '
$dCell = $ws.Range("D23")
$dCell.NumberFormat = '@'
$dCell.Value = '0.5305'
$dCell.ClearFormats()
$ws.Range("E23").Value = '  -7.40%  '
$ws.Range("E24").Value = '  -2.61%  '
$dCell = $ws.Range("D25")
$dCell.NumberFormat = '@'
$dCell.Value = '22.419.35'
$dCell.ClearFormats()
$ws.Range("E25").Value = '  -3.79%  '
$dCell = $ws.Range("D26")
$dCell.NumberFormat = '@'
$dCell.Value = '2.381'
$dCell.ClearFormats()
$ws.Range("E26").Value = '  -1.17%  '
$dCell = $ws.Range("D27")
$dCell.NumberFormat = '@'
$dCell.Value = '2.899'
$dCell.ClearFormats()
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("E28").Value = '  -3.69%  '
$dCell = $ws.Range("D29")
$dCell.NumberFormat = '@'
$dCell.Value = '146.46'
$dCell.ClearFormats()
$ws.Range("E29").Value = '  -2.80%  '
$dCell = $ws.Range("D30")
$dCell.NumberFormat = '@'
$dCell.Value = '4.984'
$dCell.ClearFormats()
$ws.Range("E30").Value = '  -3.06%  '
$ws.Range("E31").Value = '  -4.35%  '
$dCell = $ws.Range("D32")
$dCell.NumberFormat = '@'
$dCell.Value = '1.742.30'
$dCell.ClearFormats()
$dCell = $ws.Range("D33")
$dCell.NumberFormat = '@'
$dCell.Value = '1.027'
$dCell.ClearFormats()
$ws.Range("E33").Value = '  +5.34%  '
$dCell = $ws.Range("D34")
$dCell.NumberFormat = '@'
$dCell.Value = '6.240'
$dCell.ClearFormats()
$ws.Range("E34").Value = '  -7.17%  '
$dCell = $ws.Range("D35")
$dCell.NumberFormat = '@'
$dCell.Value = '2.025'
$dCell.ClearFormats()
$ws.Range("E35").Value = '  -4.34%  '
$dCell = $ws.Range("D36")
$dCell.NumberFormat = '@'
$dCell.Value = '10.14'
$dCell.ClearFormats()
$ws.Range("E36").Value = '  -9.07%  '
$dCell = $ws.Range("D37")
$dCell.NumberFormat = '@'
$dCell.Value = '0.08467'
$dCell.ClearFormats()
$ws.Range("E37").Value = '  -3.09%  '
$dCell = $ws.Range("D38")
$dCell.NumberFormat = '@'
$dCell.Value = '0.02539'
$dCell.ClearFormats()
$ws.Range("E38").Value = '  -4.34%  '
$dCell = $ws.Range("D39")
$dCell.NumberFormat = '@'
$dCell.Value = '0.2324'
$dCell.ClearFormats()
$ws.Range("E39").Value = '  -3.61%  '
$dCell = $ws.Range("D40")
$dCell.NumberFormat = '@'
$dCell.Value = '5.536'
$dCell.ClearFormats()
$ws.Range("E40").Value = '  -4.90%  '
$dCell = $ws.Range("D41")
$dCell.NumberFormat = '@'
$dCell.Value = '0.06491'
$dCell.ClearFormats()
$ws.Range("E41").Value = '  -2.98%  '
$dCell = $ws.Range("D42")
$dCell.NumberFormat = '@'
$dCell.Value = '1.304'
$dCell.ClearFormats()
$ws.Range("E42").Value = '  +1.13%  '
$dCell = $ws.Range("D43")
$dCell.NumberFormat = '@'
$dCell.Value = '11.74'
$dCell.ClearFormats()
$ws.Range("E43").Value = '  -7.13%  '
$dCell = $ws.Range("D44")
$dCell.NumberFormat = '@'
$dCell.Value = '0.6365'
$dCell.ClearFormats()
$ws.Range("E44").Value = '  -6.09%  '
$ws.Range("E45").Value = '  -6.66%  '
$ws.Range("E46").Value = '  -0.09%  '
$dCell = $ws.Range("D47")
$dCell.NumberFormat = '@'
$dCell.Value = '0.5998'
$dCell.ClearFormats()
$ws.Range("E47").Value = '  -4.68%  '
$dCell = $ws.Range("D48")
$dCell.NumberFormat = '@'
$dCell.Value = '3.762'
$dCell.ClearFormats()
$ws.Range("E48").Value = '  -3.12%  '
$dCell = $ws.Range("D49")
$dCell.NumberFormat = '@'
$dCell.Value = '2.113'
$dCell.ClearFormats()
$ws.Range("E49").Value = '  -5.03%  '
$dCell = $ws.Range("D50")
$dCell.NumberFormat = '@'
$dCell.Value = '1.262'
$dCell.ClearFormats()
$ws.Range("E50").Value = '  +3.81%  '
$dCell = $ws.Range("D51")
$dCell.NumberFormat = '@'
$dCell.Value = '125.13'
$dCell.ClearFormats()
$ws.Range("E51").Value = '  -0.70%  '
